# PinkbikeBuySellAllMountain.xlsx — remove the now-redundant
# "Attribute_Name" / "element_content_attribute_name" rows that used to
# terminate each Product_Tags block (rows 22:23 and 40:41), matching the
# companion ErrorTesting template. Row numbers below these blanks are left
# untouched (rows 25-38 keep their original row numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing "Attribute_Name"/"element_content_attribute_name"
# rows from the middle block...
$ws.Range("A22:B23").ClearContents() | Out-Null

# ...and from the final block.
$ws.Range("A40:B41").ClearContents() | Out-Null

# Leave the selection where the author ended up after deleting that block.
$ws.Range("A40:A41").Select() | Out-Null
